# Refresh LR-pair edge stats with the new TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per sending-cluster ligand stats (Ligand-expressing cells, detection rate,
# average/total expression value, derived specificity avg/total) - constant
# across all target clusters for a given sending cluster.
$sendStats = @{
  "ECs" = @(2, 0.6666666666666666, 0.265349, 0.7960469999999999, 0.1498685997319469, 0.1498685997319469)
  "FAPs" = @(3, 1, 1.003400666666667, 3.010202, 0.5667187473230926, 0.5667187473230925)
  "MuSCs" = @(3, 1, 0.4790693333333333, 1.437208, 0.2705774288246194, 0.2705774288246194)
  "Resolving-Mac" = @(1, 0.3333333333333333, 0.02272533333333333, 0.068176, 0.01283522412034115, 0.01283522412034115)
}

# Per target-cluster receptor stats (average/total expression value,
# derived specificity avg/total) - constant across all sending clusters.
$targetStats = @{
  "ECs" = @(13.67700833333333, 41.031025, 0.124413831206147, 0.124413831206147)
  "FAPs" = @(74.64939600000001, 223.948188, 0.679053278848249, 0.6790532788482488)
  "Inflammatory-Mac" = @(1.629335666666667, 4.888007, 0.01482136207497777, 0.01482136207497777)
  "MuSCs" = @(19.17462033333333, 57.523861, 0.174423230537864, 0.174423230537864)
  "Resolving-Mac" = @(0.801214, 2.403642, 0.007288297332762355, 0.007288297332762355)
}

# Per (sending, target) edge weight / derived specificity (avg & total).
$edgeStats = @{
  "ECs|ECs" = @(3.629180484241667, 32.66262435817499, 0.01864572667015204, 0.01864572667015204)
  "ECs|FAPs" = @(19.808142579204, 178.273283212836, 0.1017687640443743, 0.1017687640443743)
  "ECs|Inflammatory-Mac" = @(0.4323425898143334, 3.891083308329, 0.002221256780297101, 0.002221256780297101)
  "ECs|MuSCs" = @(5.087966330829667, 45.79169697746699, 0.02614056532143223, 0.02614056532143223)
  "ECs|Resolving-Mac" = @(0.212601333686, 1.913412003174, 0.001092286915691177, 0.001092286915691177)
  "FAPs|ECs" = @(13.72351927967222, 123.51167351705, 0.07050765057081429, 0.07050765057081428)
  "FAPs|FAPs" = @(74.90325371266401, 674.1292834139759, 0.3848322235545183, 0.3848322235545182)
  "FAPs|Inflammatory-Mac" = @(1.634876494157111, 14.713888447414, 0.008399543748753395, 0.008399543748753393)
  "FAPs|MuSCs" = @(19.23982682554689, 173.158441429922, 0.09884891471446526, 0.09884891471446525)
  "FAPs|Resolving-Mac" = @(0.8039386617426666, 7.235447955683999, 0.004130414734541319, 0.004130414734541318)
  "MuSCs|ECs" = @(6.552235264244444, 58.97011737819999, 0.03366357455797945, 0.03366357455797945)
  "MuSCs|FAPs" = @(35.762236375456, 321.860127379104, 0.1837364902256865, 0.1837364902256865)
  "MuSCs|Inflammatory-Mac" = @(0.7805647516062222, 7.025082764455999, 0.004010326041926212, 0.004010326041926211)
  "MuSCs|MuSCs" = @(9.185972580009777, 82.67375322008799, 0.04719498924621908, 0.04719498924621908)
  "MuSCs|Resolving-Mac" = @(0.3838370568373333, 3.454533511536, 0.00197204875280817, 0.00197204875280817)
  "Resolving-Mac|ECs" = @(0.3108145733777778, 2.7973311604, 0.001596879407201189, 0.001596879407201189)
  "Resolving-Mac|FAPs" = @(1.696432407232, 15.267891665088, 0.008715801023669787, 0.008715801023669785)
  "Resolving-Mac|Inflammatory-Mac" = @(0.03702719613688889, 0.333244765232, 0.0001902355040010642, 0.0001902355040010642)
  "Resolving-Mac|MuSCs" = @(0.4357496386151111, 3.921746747536, 0.002238761255747416, 0.002238761255747416)
  "Resolving-Mac|Resolving-Mac" = @(0.01820785522133333, 0.163870696992, 0.00009354692972168942, 0.00009354692972168942)
}

$sendOrder = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")
$targetOrder = @("ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Resolving-Mac")

$row = 2
foreach ($send in $sendOrder) {
  $ss = $sendStats[$send]
  foreach ($target in $targetOrder) {
    $ts = $targetStats[$target]
    $es = $edgeStats["$send|$target"]
    $ws.Cells.Item($row, 1).Value = $send
    $ws.Cells.Item($row, 2).Value = "Fgf1"
    $ws.Cells.Item($row, 3).Value = "Fgfr1"
    $ws.Cells.Item($row, 4).Value = $target
    $ws.Cells.Item($row, 5).Value = $ss[0]
    $ws.Cells.Item($row, 6).Value = $ss[1]
    $ws.Cells.Item($row, 7).Value = $ss[2]
    $ws.Cells.Item($row, 8).Value = $ss[3]
    $ws.Cells.Item($row, 9).Value = $ss[4]
    $ws.Cells.Item($row, 10).Value = $ss[5]
    $ws.Cells.Item($row, 11).Value = 3
    $ws.Cells.Item($row, 12).Value = 1
    $ws.Cells.Item($row, 13).Value = $ts[0]
    $ws.Cells.Item($row, 14).Value = $ts[1]
    $ws.Cells.Item($row, 15).Value = $ts[2]
    $ws.Cells.Item($row, 16).Value = $ts[3]
    $ws.Cells.Item($row, 17).Value = $es[0]
    $ws.Cells.Item($row, 18).Value = $es[1]
    $ws.Cells.Item($row, 19).Value = $es[2]
    $ws.Cells.Item($row, 20).Value = $es[3]
    $row++
  }
}
